$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the "Status" values here are driven from the same shared
# text as the per-language sheets, so keep them in sync with the new wording.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: record the handback for both tracked source files.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B2").Value = $statusHandedBack
$wsZhCn.Range("E2").Value = "36348082-362f-4f10-8e2e-cf196d7b98a8.md"
$wsZhCn.Range("F2").Value = "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-01-26 04:33:05"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("E3").Value = "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md"
$wsZhCn.Range("F3").Value = "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "2016-01-26 04:33:05"
$wsZhCn.Range("H3").Value = "Include"

# Rebuild the hyperlinks in file order (A2,C2,E2,F2,A3,C3,E3,F3,A4) so the
# relationship ids line up the same way Excel would renumber them.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/36348082-362f-4f10-8e2e-cf196d7b98a8.md", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5344bf262e4ac21c4dfdee9e57ba99756893a7c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.zh-cn.xlf", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/36348082-362f-4f10-8e2e-cf196d7b98a8.md", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5344bf262e4ac21c4dfdee9e57ba99756893a7c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.zh-cn.xlf", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/caf0c8ef-ce87-4ce9-bde5-72454b17a936.md", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5344bf262e4ac21c4dfdee9e57ba99756893a7c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.zh-cn.xlf", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/caf0c8ef-ce87-4ce9-bde5-72454b17a936.md", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5344bf262e4ac21c4dfdee9e57ba99756893a7c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.zh-cn.xlf", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet: same handback bookkeeping as zh-cn, different target files.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B2").Value = $statusHandedBack
$wsDeDe.Range("E2").Value = "36348082-362f-4f10-8e2e-cf196d7b98a8.md"
$wsDeDe.Range("F2").Value = "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.de-de.xlf"
$wsDeDe.Range("G2").Value = "2016-01-26 04:33:18"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("E3").Value = "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md"
$wsDeDe.Range("F3").Value = "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.de-de.xlf"
$wsDeDe.Range("G3").Value = "2016-01-26 04:33:18"
$wsDeDe.Range("H3").Value = "Include"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/36348082-362f-4f10-8e2e-cf196d7b98a8.md", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd852d56b9bd6c992318cb33068f554a2bd2539c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.de-de.xlf", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/36348082-362f-4f10-8e2e-cf196d7b98a8.md", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd852d56b9bd6c992318cb33068f554a2bd2539c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.de-de.xlf", "", "", "36348082-362f-4f10-8e2e-cf196d7b98a8.1a94036de4c812580f2dbd31ea60eb55e630ea29.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/caf0c8ef-ce87-4ce9-bde5-72454b17a936.md", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd852d56b9bd6c992318cb33068f554a2bd2539c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.de-de.xlf", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/e2e/caf0c8ef-ce87-4ce9-bde5-72454b17a936.md", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd852d56b9bd6c992318cb33068f554a2bd2539c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.de-de.xlf", "", "", "caf0c8ef-ce87-4ce9-bde5-72454b17a936.a1e707a14eb50eafab9f10cdbdd76cd30a79a09e.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/56f2b0c1f01f921bb2f4d4acf7090a2730389b6c/.localization-config", "", "", ".localization-config")
